$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (Through 2022-07-17 -> Through 2022-07-18)
$ws.Name = "Through 2022-07-18"

# Update the header label in I1 (shared string "2022 (through 07-17)" -> "2022 (through 07-18)")
$ws.Range("I1").Value = "2022 (through 07-18)"

# Update data values for new day (2022-07-26 data add)
$ws.Range("I8").Value = 99
$ws.Range("I14").Value = 904
